# Add StdOut to Runcommand
#
# Duplicates the "dbo#ExcelTest" worksheet twice (Excel's native "Move or
# Copy" behaviour), inserting the copies between "dbo#ExcelTest" and
# "Assert" and naming them "dbo#ExcelTest2" / "dbo#ExcelTest3". The last
# copy ("dbo#ExcelTest3") gets a tweaked E10 value (text "4" instead of the
# numeric 4) and becomes the active sheet/selection.

$wb = $excel.ActiveWorkbook

$base = $wb.Worksheets.Item("dbo#ExcelTest")

# First copy -> dbo#ExcelTest2, placed right after dbo#ExcelTest.
$base.Copy($null, $base)
$copy1 = $wb.Worksheets.Item(2)
$copy1.Name = "dbo#ExcelTest2"

# Second copy (of the first copy) -> dbo#ExcelTest3, placed right after
# dbo#ExcelTest2 (still before Assert).
$copy1.Copy($null, $copy1)
$copy2 = $wb.Worksheets.Item(3)
$copy2.Name = "dbo#ExcelTest3"

# On the new third sheet, row 10's Key column holds the text "4" rather
# than the numeric 4 used on the other copies.
$copy2.Range("E10").Value = "4"

# Leave the selection/active sheet on dbo#ExcelTest3, cell E11.
$copy2.Range("E11").Select() | Out-Null
